# Fill in the "Result" column (S) with "Pass" for every row whose
# RunMode (column A) is "Yes" on the TestData sheet, then leave the
# selection / view the way Excel recorded it after this edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $runMode = $ws.Cells.Item($r, 1).Value2
    if ($runMode -eq "Yes") {
        $ws.Cells.Item($r, 19).Value = "Pass"   # column S = 19
    }
}

# Match the recorded view/selection state from the saved workbook.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 10
$ws.Range("S2:S10").Select()
